# Update the "dSF" column (F) values for the wicks_jordan sheet.
# These new values reflect a repull/recalculation of the underlying data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -4
    4  = -3
    5  = 3
    7  = -1
    8  = 3
    10 = -4
    12 = 0
    13 = 1
    14 = -1
    15 = 4
    16 = -1
    18 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
